# LSR1_H_SOE.xlsx - RmD and SoE updates
# Inserts a new "Timepoint" column into the SOE summary-of-evidence table
# and re-flows the column widths / row heights / view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column B ("Timepoint") ------------------------------
$ws.Columns("B:B").Insert()

# --- Header row ---------------------------------------------------------
$ws.Range("B1").Value = "Timepoint"

# --- Timepoint values for each evidence row -----------------------------
$ws.Range("B2").Value = "4-12 weeks"
$ws.Range("B3").Value = "4-12 weeks"
$ws.Range("B4").Value = "8 weeks"
$ws.Range("B5").Value = "8 weeks"
$ws.Range("B6").Value = "8 weeks"
$ws.Range("B7").Value = "8 weeks"
$ws.Range("B8").Value = "8 weeks"
$ws.Range("B9").Value = "8 weeks"

# --- Column widths --------------------------------------------------------
$ws.Columns("A").ColumnWidth = 16.90625
$ws.Columns("B").ColumnWidth = 16.90625
$ws.Columns("C").ColumnWidth = 20.36328125
$ws.Columns("D").ColumnWidth = 24.453125
$ws.Columns("E").ColumnWidth = 21
$ws.Columns("F").ColumnWidth = 25.453125
$ws.Columns("G").ColumnWidth = 23.81640625

# --- Row heights (re-flowed after widening data columns) -----------------
$ws.Rows(1).RowHeight = 40.5
$ws.Rows(2).RowHeight = 189
$ws.Rows(3).RowHeight = 108
$ws.Rows(4).RowHeight = 94.5
$ws.Rows(10).RowHeight = 81.5
$ws.Rows(11).RowHeight = 94.5

# --- View state -------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("E9").Select()
